# Insert a new data row at row 40 (weekly price-report update).
# This pushes the existing rows 40..128 down to 41..129 and keeps the
# sheet dimension / shared data intact (Excel handles that automatically
# via Rows().Insert(), matching a manual "Insert Row" in the UI).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(40).Insert()

# Populate the freshly inserted row 40 with the new record.
$ws.Cells.Item(40, 1).Value  = 3
$ws.Cells.Item(40, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(40, 3).Value  = "Coquimbo"
$ws.Cells.Item(40, 4).Value  = 44581
$ws.Cells.Item(40, 5).Value  = 5
$ws.Cells.Item(40, 6).Value  = 100112030
$ws.Cells.Item(40, 7).Value  = "Poroto granado"
$ws.Cells.Item(40, 8).Value  = "Sin especificar"
$ws.Cells.Item(40, 9).Value  = "Primera"
$ws.Cells.Item(40, 10).Value = 73
$ws.Cells.Item(40, 11).Value = 25000
$ws.Cells.Item(40, 12).Value = 26000
$ws.Cells.Item(40, 13).Value = 25479
$ws.Cells.Item(40, 14).Value = '$/malla 25 kilos'
$ws.Cells.Item(40, 15).Value = "Provincia de Talca"
$ws.Cells.Item(40, 16).Value = 1019
$ws.Cells.Item(40, 17).Value = 25
$ws.Cells.Item(40, 18).Value = "Hortaliza"
